# Weekly update: insert a new daily price record as the new row 4,
# pushing the existing rows (previously 4-17) down to 5-18.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4 (shifts rows 4..17 down to 5..18).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new observation.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44537
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 800
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = 850
$ws.Range("N4").Value = "`$/kilo"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 850
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
